$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current "Yco" column (H) to hold the new
# "CO2/(CO+CO2)" ratio column. This shifts the existing H:O columns to I:P.
$ws.Range("H1").EntireColumn.Insert()

# New column keeps the same visual width as the columns to its left
# (D:G, 14.33203125 raw width <-> ~13.5 character width).
$ws.Range("H1").ColumnWidth = 13.5

# Fill in the formula for the new ratio column: feedYco2 / (feedYco + feedYco2)
$ws.Range("H2").Formula = "=F2/(E2+F2)"
$ws.Range("H3:H30").Formula = "=F3/(E3+F3)"

# Give the new column the same (bordered, header-like) style as the header
# row uses, matching the rest of the data columns' look.
$ws.Range("H1").Copy()
$ws.Range("H2:H30").PasteSpecial(-4122)

# Header text for the new column.
$ws.Range("H1").Value = "CO2/(CO+CO2)"

# Match the saved selection in the source workbook.
$ws.Range("H2").Select()
